$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for data rows 2-77 from 45178 (2023-09-09)
# to 45179 (2023-09-10), matching the automatic update reflected in the diff.
for ($row = 2; $row -le 77; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}
